$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60 previously only had an (empty) styled cell in column C. Fill in the
# new LeetCode entry "184. Department Highest Salary" (a SQL problem).
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "184. Department Highest Salary"
$ws.Range("C60").Value = "https://leetcode.com/problems/department-highest-salary/"
$ws.Range("D60").Value = "SQL"
$ws.Range("E60").Value = "Medium"
$ws.Range("F60").Value = "Solved"

# Move the selection to reflect where the author ended up after adding the row.
$ws.Range("C66").Select() | Out-Null
